$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentários")

# Link the blank input cells in column B to the corresponding confusion-matrix
# counts on the 'Matriz Confusão 2' sheet instead of leaving them manually blank.
$ws.Range("B2").Formula = "='Matriz Confusão 2'!C3"
$ws.Range("B3").Formula = "='Matriz Confusão 2'!D3"
$ws.Range("B4").Formula = "='Matriz Confusão 2'!E3"
$ws.Range("B6").Formula = "='Matriz Confusão 2'!D4"
$ws.Range("B7").Formula = "='Matriz Confusão 2'!C4"
$ws.Range("B8").Formula = "='Matriz Confusão 2'!E4"
$ws.Range("B10").Formula = "='Matriz Confusão 2'!E5"
$ws.Range("B11").Formula = "='Matriz Confusão 2'!C5"
$ws.Range("B12").Formula = "='Matriz Confusão 2'!D5"

# Update the last selected cell on the "Comentários" sheet.
$ws.Range("I9").Select()

# Restore "Geral" as the active sheet/tab, matching the original workbook state.
$wb.Worksheets.Item("Geral").Activate()
